$d = $word.ActiveDocument

# --- Paragraph 6: "If(EX/MEM.RegWrite AND (EX/MEM.Rd == ID/EX.Rs)) " ---
# -> "If(EX/MEM.RegWrite AND (EX/MEM.Rd == ID/EX.(Rs/Rd))) "
$r = $d.Paragraphs(6).Range
$r.Find.Execute("ID/EX.Rs", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "ID/EX.(Rs/Rd)", 2)

# --- Paragraph 7: "Forward EX/MEM.Rd data" -> "Forward EX/MEM.Rd data " ---
$p7 = $d.Paragraphs(7).Range
$ins = $d.Range($p7.Start, $p7.End - 1)
$ins.Collapse(0)
$ins.InsertAfter(" ")

# --- Paragraph 10: "If(MEM/WP.RegWrite AND (MEM/WP.Rd == ID/EX.Rs)) " ---
# -> "If(MEM/WP.RegWrite AND (MEM/WP.Rd == ID/EX.(Rs/Rd))) "
$r = $d.Paragraphs(10).Range
$r.Find.Execute("ID/EX.Rs", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "ID/EX.(Rs/Rd)", 2)

# --- Paragraph 11: "Forward MEM/WP.Rd data" -> "Forward MEM/WP.Rd data to " ---
$p11 = $d.Paragraphs(11).Range
$ins = $d.Range($p11.Start, $p11.End - 1)
$ins.Collapse(0)
$ins.InsertAfter(" to ")

# --- Paragraph 13: "If(MEM/WP.RegWrite AND NOT(EX/MEM.RegWrite) AND (EX/MEM.Rd != ID/EX.Rs)) " ---
# -> "If(MEM/WP.RegWrite AND NOT(EX/MEM.RegWrite) AND (EX/MEM.Rd != ID/EX.(Rs/Rd))) "
$r = $d.Paragraphs(13).Range
$r.Find.Execute("ID/EX.Rs", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "ID/EX.(Rs/Rd)", 2)

# --- Paragraph 16: "If(ID/EX.MemRead AND (ID/EX.RegisterRs = IF/ID.RegisterRs))" ---
# -> "If(ID/EX.MemRead AND (ID/EX.Rs = IF/ID.(Rs/Rd)))"
$r = $d.Paragraphs(16).Range
$r.Find.Execute("RegisterRs", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "Rs", 2)
$r = $d.Paragraphs(16).Range
$r.Find.Execute("ID.Rs", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "ID.(Rs/Rd)", 2)

# Move the _GoBack bookmark from paragraph 6 to paragraph 16, right after "If(ID/EX.MemRead AND (ID/EX."
$p16 = $d.Paragraphs(16).Range
$pos = $p16.Start + 28
$bm = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bm)
